# Auto-applies the cell-level market-data refresh captured in the commit diff.
# Each worksheet is an Excel Table (Table_<CLASS>) of FFXIV Leve profit calculations;
# the scheduled runner only overwrites the raw price/profit columns (H:N) per row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1728.2
$ws.Range("I92").Value = 1789.8462
$ws.Range("J92").Value = 1327.5
$ws.Range("K92").Value = 1789.8462
$ws.Range("L92").Value = 1327.5
$ws.Range("M92").Value = -541.8462
$ws.Range("N92").Value = -3823.5

$ws.Range("H98").Value = 1451.7858
$ws.Range("I98").Value = 1451.7858
$ws.Range("K98").Value = 1451.7858
$ws.Range("M98").Value = 46.21419999999989

$ws.Range("H122").Value = 1451.7858
$ws.Range("I122").Value = 1451.7858
$ws.Range("K122").Value = 4355.357400000001
$ws.Range("M122").Value = -1905.357400000001

$ws.Range("H137").Value = 9435283
$ws.Range("I137").Value = 1302.8864
$ws.Range("J137").Value = 55556964
$ws.Range("K137").Value = 3908.6592
$ws.Range("L137").Value = 166670892
$ws.Range("M137").Value = -1358.6592
$ws.Range("N137").Value = -166675992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 852.0465
$ws.Range("I74").Value = 858.06665
$ws.Range("J74").Value = 838.1539
$ws.Range("K74").Value = 858.06665
$ws.Range("L74").Value = 838.1539
$ws.Range("M74").Value = 15.93335000000002
$ws.Range("N74").Value = -2586.1539

$ws.Range("H77").Value = 852.0465
$ws.Range("I77").Value = 858.06665
$ws.Range("J77").Value = 838.1539
$ws.Range("K77").Value = 4290.33325
$ws.Range("L77").Value = 4190.7695
$ws.Range("M77").Value = 77.66675000000032
$ws.Range("N77").Value = -12926.7695

$ws.Range("H132").Value = 1381.5968
$ws.Range("I132").Value = 937.4681
$ws.Range("J132").Value = 2773.2
$ws.Range("K132").Value = 2812.4043
$ws.Range("L132").Value = 8319.599999999999
$ws.Range("M132").Value = -282.4043000000001
$ws.Range("N132").Value = -13379.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 2004
$ws.Range("I8").Value = 2004
$ws.Range("K8").Value = 2004
$ws.Range("M8").Value = -1864

$ws.Range("H20").Value = 50033536
$ws.Range("I20").Value = 216502.67
$ws.Range("K20").Value = 216502.67
$ws.Range("M20").Value = -216255.67

$ws.Range("H27").Value = 17500
$ws.Range("J27").Value = 17500
$ws.Range("L27").Value = 17500
$ws.Range("N27").Value = -17884

$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 20618.182
$ws.Range("J68").Value = 20618.182
$ws.Range("L68").Value = 20618.182
$ws.Range("N68").Value = -22116.182

$ws.Range("H71").Value = 20618.182
$ws.Range("J71").Value = 20618.182
$ws.Range("L71").Value = 61854.546
$ws.Range("N71").Value = -69342.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 155.05882
$ws.Range("I6").Value = 139.41667
$ws.Range("J6").Value = 192.6
$ws.Range("K6").Value = 418.25001
$ws.Range("L6").Value = 577.8
$ws.Range("M6").Value = -305.25001
$ws.Range("N6").Value = -803.8

$ws.Range("H70").Value = 2916.875
$ws.Range("I70").Value = 2805.8333
$ws.Range("J70").Value = 3250
$ws.Range("K70").Value = 8417.499899999999
$ws.Range("L70").Value = 9750
$ws.Range("M70").Value = -8102.499899999999
$ws.Range("N70").Value = -10380

$ws.Range("H73").Value = 2916.875
$ws.Range("I73").Value = 2805.8333
$ws.Range("J73").Value = 3250
$ws.Range("K73").Value = 8417.499899999999
$ws.Range("L73").Value = 9750
$ws.Range("M73").Value = -7325.499899999999
$ws.Range("N73").Value = -11934

$ws.Range("H75").Value = 431.6
$ws.Range("J75").Value = 381.66666
$ws.Range("L75").Value = 1144.99998
$ws.Range("N75").Value = -3140.99998

$ws.Range("H78").Value = 431.6
$ws.Range("J78").Value = 381.66666
$ws.Range("L78").Value = 3434.99994
$ws.Range("N78").Value = -13418.99994

$ws.Range("H87").Value = 16163
$ws.Range("I87").Value = 8399.200000000001
$ws.Range("J87").Value = 21708.572
$ws.Range("K87").Value = 25197.6
$ws.Range("L87").Value = 65125.716
$ws.Range("M87").Value = -23949.6
$ws.Range("N87").Value = -67621.716

$ws.Range("H90").Value = 16163
$ws.Range("I90").Value = 8399.200000000001
$ws.Range("J90").Value = 21708.572
$ws.Range("K90").Value = 75592.8
$ws.Range("L90").Value = 195377.148
$ws.Range("M90").Value = -69352.8
$ws.Range("N90").Value = -207857.148

$ws.Range("H103").Value = 548.6
$ws.Range("I103").Value = 435.75
$ws.Range("J103").Value = 1000
$ws.Range("K103").Value = 1307.25
$ws.Range("L103").Value = 3000
$ws.Range("M103").Value = -428.25
$ws.Range("N103").Value = -4758

$ws.Range("H131").Value = 44884.668
$ws.Range("I131").Value = 101513
$ws.Range("J131").Value = 4435.857
$ws.Range("K131").Value = 304539
$ws.Range("L131").Value = 13307.571
$ws.Range("M131").Value = -299499
$ws.Range("N131").Value = -23387.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4423.6206
$ws.Range("I70").Value = 4067.5
$ws.Range("J70").Value = 5542.857
$ws.Range("K70").Value = 4067.5
$ws.Range("L70").Value = 5542.857
$ws.Range("M70").Value = -3797.5
$ws.Range("N70").Value = -6082.857

$ws.Range("H73").Value = 4423.6206
$ws.Range("I73").Value = 4067.5
$ws.Range("J73").Value = 5542.857
$ws.Range("K73").Value = 4067.5
$ws.Range("L73").Value = 5542.857
$ws.Range("M73").Value = -3131.5
$ws.Range("N73").Value = -7414.857

$ws.Range("H119").Value = 25666.666
$ws.Range("J119").Value = 25666.666
$ws.Range("L119").Value = 25666.666
$ws.Range("N119").Value = -35342.666

$ws.Range("H122").Value = 44044424
$ws.Range("I122").Value = 59588812
$ws.Range("J122").Value = 1984.6666
$ws.Range("K122").Value = 178766436
$ws.Range("L122").Value = 5953.9998
$ws.Range("M122").Value = -178763986
$ws.Range("N122").Value = -10853.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1169.76
$ws.Range("I22").Value = 215.8
$ws.Range("J22").Value = 1408.25
$ws.Range("K22").Value = 215.8
$ws.Range("L22").Value = 1408.25
$ws.Range("M22").Value = 79.19999999999999
$ws.Range("N22").Value = -1998.25

$ws.Range("H27").Value = 1169.76
$ws.Range("I27").Value = 215.8
$ws.Range("J27").Value = 1408.25
$ws.Range("K27").Value = 215.8
$ws.Range("L27").Value = 1408.25
$ws.Range("M27").Value = -108.8
$ws.Range("N27").Value = -1622.25

$ws.Range("H40").Value = 13751046
$ws.Range("I40").Value = 14323840
$ws.Range("K40").Value = 14323840
$ws.Range("M40").Value = -14323704

$ws.Range("H139").Value = 49543
$ws.Range("J139").Value = 49543
$ws.Range("L139").Value = 49543
$ws.Range("N139").Value = -59823

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 27642
$ws.Range("J93").Value = 27642
$ws.Range("L93").Value = 27642
$ws.Range("N93").Value = -32634
